# Tied DCIN to VSYS, added LED, cleaned up Bridge Rectifiers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1. Remove the old "DCIN tied to VSYS" resistor row (Resistor (0ohm) /
#    R-US_R1210) that used to live at row 6. Deleting the row shifts
#    every row below it up by one.
# ------------------------------------------------------------------
$ws.Rows("6:6").Delete()

# ------------------------------------------------------------------
# 2. After the shift:
#    row 11 -> Capacitor (22uF)   needs a part code + new datasheet link
#    row 12 -> Capacitor (47pF)   needs a part code (same datasheet link)
#    row 15 -> Inductor (1uH)     needs a part code + new datasheet link
# ------------------------------------------------------------------
$ws.Range("D11").Value = "GRM188R60G226MEA0D"
$ws.Range("D12").Value = "CL10C470JB8NNNC"
$ws.Range("D15").Value = "LQM32PN1R0MG0L"

# Point the existing Capacitor (22uF) hyperlink at the new datasheet.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$G$11') {
        $h.Address = "http://www.digikey.com/product-detail/en/GRM188R60G226MEA0D/490-5526-1-ND/2334922"
    }
}

# Point the existing Inductor (1uH) hyperlink at the new datasheet.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$G$15') {
        $h.Address = "http://www.digikey.com/product-detail/en/LQM32PN1R0MG0L/490-10767-1-ND/5251332"
    }
}

# ------------------------------------------------------------------
# 3. Fill in the two previously-blank spacer rows (16 and 17) with the
#    new Resistor (1K) and LED parts.
# ------------------------------------------------------------------
$ws.Range("C16").Value = "Resistor (1K)"
$ws.Range("D16").Value = "RC0603JR-071KL"
$ws.Range("E16").Value = "RES-07856"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = "http://www.digikey.com/product-detail/en/RC0603JR-071KL/311-1.0KGRCT-ND/729624"
$ws.Hyperlinks.Add($ws.Range("G16"), "http://www.digikey.com/product-detail/en/RC0603JR-071KL/311-1.0KGRCT-ND/729624")
$ws.Range("G16").Style = "Hyperlink"

$ws.Range("C17").Value = "LED"
$ws.Range("D17").Value = "LTST-C191TBKT"
$ws.Range("E17").Value = "DIO-08575"
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = "http://www.digikey.com/product-detail/en/LTST-C191TBKT/160-1647-1-ND/573587"
$ws.Hyperlinks.Add($ws.Range("G17"), "http://www.digikey.com/product-detail/en/LTST-C191TBKT/160-1647-1-ND/573587")
$ws.Range("G17").Style = "Hyperlink"

# ------------------------------------------------------------------
# 4. Column D got a bit wider to fit the longer part codes.
# ------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 20.83

# ------------------------------------------------------------------
# 5. Refresh the view: active cell / scroll position.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("D20").Select()
